$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing data (rows 1-25) stores every value - including things that
# look numeric, like "3" or "66.67%" - as literal text. Force the new cells
# to the same "Text" storage so Excel doesn't auto-coerce them into numbers
# or percentages, then reset the display style back to Normal/General so the
# new rows don't end up with a different cell style than the rest of the
# sheet.
$newRange = $ws.Range("A26:I27")
$newRange.NumberFormat = "@"

$ws.Range("A26").Value = "pumas unam"
$ws.Range("B26").Value = "3"
$ws.Range("C26").Value = "2"
$ws.Range("D26").Value = "1"
$ws.Range("E26").Value = "1"
$ws.Range("F26").Value = "66.67%"
$ws.Range("G26").Value = "66.67%"
$ws.Range("H26").Value = "66.67%"
$ws.Range("I26").Value = "50.00%"

$ws.Range("A27").Value = "pachuca"
$ws.Range("B27").Value = "3"
$ws.Range("C27").Value = "1"
$ws.Range("D27").Value = "1"
$ws.Range("E27").Value = "1"
$ws.Range("F27").Value = "33.33%"
$ws.Range("G27").Value = "33.33%"
$ws.Range("H27").Value = "33.33%"
$ws.Range("I27").Value = "50.00%"

$newRange.Style = "Normal"
